$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text runs (Volume/Number and date range) ---
# "Volume 31   Number  46" -> "...47"
$ws.Range("A8").Characters(21,2).Text = "47"
# "Report Covering the Week  11/11/2024  Through  11/17/2024"
#  -> week shifted forward by one (11/18/2024 .. 11/24/2024)
$ws.Range("C9").Characters(27,10).Text = "11/18/2024"
$ws.Range("C9").Characters(48,10).Text = "11/24/2024"

# --- Simple numeric value updates (style/type unchanged) ---
$ws.Range("H15").Value = -100
$ws.Range("M15").Value = -11.111111111111
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 200
$ws.Range("I16").Value = 109
$ws.Range("K16").Value = 31.325301204819
$ws.Range("L16").Value = 1.869158878504
$ws.Range("M16").Value = 39.743589743589
$ws.Range("N16").Value = -79.963235294117
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 5
$ws.Range("H17").Value = 160
$ws.Range("I17").Value = 111
$ws.Range("J17").Value = 97
$ws.Range("K17").Value = 14.432989690721
$ws.Range("L17").Value = 16.842105263157
$ws.Range("M17").Value = 109.433962264151
$ws.Range("N17").Value = 26.136363636363
$ws.Range("C18").Value = 6
$ws.Range("F18").Value = 11
$ws.Range("H18").Value = 120
$ws.Range("I18").Value = 128
$ws.Range("J18").Value = 98
$ws.Range("K18").Value = 30.612244897959
$ws.Range("L18").Value = 3.225806451612
$ws.Range("M18").Value = 30.612244897959
$ws.Range("N18").Value = -84.559710494571
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -12.5
$ws.Range("F19").Value = 58
$ws.Range("G19").Value = 68
$ws.Range("H19").Value = -14.705882352941
$ws.Range("I19").Value = 707
$ws.Range("J19").Value = 668
$ws.Range("K19").Value = 5.838323353293
$ws.Range("L19").Value = -5.858854860186
$ws.Range("M19").Value = 15.711947626841
$ws.Range("N19").Value = -56.94275274056
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 40
$ws.Range("I20").Value = 53
$ws.Range("J20").Value = 90
$ws.Range("K20").Value = -41.111111111111
$ws.Range("L20").Value = -22.058823529411
$ws.Range("M20").Value = 76.666666666666
$ws.Range("N20").Value = -94.95238095238
$ws.Range("C21").Value = 26
$ws.Range("E21").Value = 30
$ws.Range("F21").Value = 104
$ws.Range("G21").Value = 89
$ws.Range("H21").Value = 16.853932584269
$ws.Range("I21").Value = 1116
$ws.Range("J21").Value = 1050
$ws.Range("K21").Value = 6.285714285714
$ws.Range("L21").Value = -3.543647363872
$ws.Range("M21").Value = 26.818181818181
$ws.Range("N21").Value = -73.263057019645
$ws.Range("F22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("J23").Value = 35
$ws.Range("K23").Value = 17.142857142857
$ws.Range("C24").Value = 34
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 88.888888888888
$ws.Range("F24").Value = 138
$ws.Range("G24").Value = 83
$ws.Range("H24").Value = 66.265060240963
$ws.Range("I24").Value = 1229
$ws.Range("J24").Value = 1053
$ws.Range("K24").Value = 16.714150047483
$ws.Range("L24").Value = 8.281938325991
$ws.Range("M24").Value = 28.422152560083
$ws.Range("C25").Value = 29
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = 93.333333333333
$ws.Range("F25").Value = 115
$ws.Range("G25").Value = 79
$ws.Range("H25").Value = 45.569620253164
$ws.Range("I25").Value = 994
$ws.Range("J25").Value = 818
$ws.Range("K25").Value = 21.515892420537
$ws.Range("L25").Value = 7.575757575757
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 40
$ws.Range("F26").Value = 23
$ws.Range("G26").Value = 11
$ws.Range("H26").Value = 109.090909090909
$ws.Range("I26").Value = 214
$ws.Range("J26").Value = 204
$ws.Range("K26").Value = 4.901960784313
$ws.Range("L26").Value = 4.901960784313
$ws.Range("M26").Value = -15.748031496063
$ws.Range("H27").Value = -100
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 2
$ws.Range("F28").Value = 8
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 33.333333333333
$ws.Range("I28").Value = 36
$ws.Range("J28").Value = 38
$ws.Range("K28").Value = -5.263157894736
$ws.Range("L28").Value = -25
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = -50

# --- Number -> placeholder text ("0" / "***.*") conversions ---
# Force text storage via NumberFormat "@" so the numeric-looking
# placeholder strings are not re-parsed back into numbers, then
# restore the workbook-standard placeholder style (copied from A15,
# which is already styled as a text placeholder cell) so the saved
# style index matches what the rest of the sheet uses (s="13").
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "0"
$ws.Range("A15").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0"
$ws.Range("A15").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "***.*"
$ws.Range("A15").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("A15").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("A15").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "0"
$ws.Range("A15").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("F27").NumberFormat = "@"
$ws.Range("F27").Value = "0"
$ws.Range("A15").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0"
$ws.Range("A15").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "***.*"
$ws.Range("A15").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "0"
$ws.Range("A15").Copy()
$ws.Range("G33").PasteSpecial(-4122)
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "***.*"
$ws.Range("A15").Copy()
$ws.Range("H33").PasteSpecial(-4122)

# --- Text placeholder -> number conversions ---
# Restore the column-appropriate numeric style afterwards: F18 carries
# the standard "count" number style, H18 the standard "percent" style.
$ws.Range("D18").Value = 2
$ws.Range("F18").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E18").Value = 200
$ws.Range("H18").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("C20").Value = 1
$ws.Range("F18").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("D20").Value = 1
$ws.Range("F18").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").Value = 0
$ws.Range("H18").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("D23").Value = 1
$ws.Range("F18").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").Value = -100
$ws.Range("H18").Copy()
$ws.Range("E23").PasteSpecial(-4122)

$excel.CutCopyMode = $false